# Updates the FlashScore odds/correct-score table (rows 4-25) to the
# values published in the 2024-12-05 refresh. Cells are addressed via
# Cells.Item(row, columnIndex) to avoid A1-notation typos across the
# ~50 "Odd_*" columns (G..BD).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4, 7).Value = 3.6
$ws.Cells.Item(4, 8).Value = 3.5
$ws.Cells.Item(4, 9).Value = 2
$ws.Cells.Item(4, 10).Value = 4
$ws.Cells.Item(4, 11).Value = 2.2
$ws.Cells.Item(4, 12).Value = 2.63
$ws.Cells.Item(4, 13).Value = 1.05
$ws.Cells.Item(4, 14).Value = 11
$ws.Cells.Item(4, 17).Value = 1.83
$ws.Cells.Item(4, 18).Value = 1.98
$ws.Cells.Item(4, 21).Value = 1.67
$ws.Cells.Item(4, 22).Value = 2.1
$ws.Cells.Item(4, 23).Value = 12
$ws.Cells.Item(4, 24).Value = 19
$ws.Cells.Item(4, 29).Value = 11
$ws.Cells.Item(4, 30).Value = 6.5
$ws.Cells.Item(4, 35).Value = 10
$ws.Cells.Item(4, 36).Value = 9
$ws.Cells.Item(4, 40).Value = 5.5
$ws.Cells.Item(4, 41).Value = 19
$ws.Cells.Item(4, 47).Value = 7.5
$ws.Cells.Item(4, 51).Value = 11
$ws.Cells.Item(4, 52).Value = 21

$ws.Cells.Item(5, 7).Value = 2.15
$ws.Cells.Item(5, 8).Value = 3.3
$ws.Cells.Item(5, 9).Value = 3.5
$ws.Cells.Item(5, 26).Value = 19
$ws.Cells.Item(5, 27).Value = 17
$ws.Cells.Item(5, 37).Value = 41
$ws.Cells.Item(5, 39).Value = 41
$ws.Cells.Item(5, 40).Value = 4
$ws.Cells.Item(5, 50).Value = 5.5

$ws.Cells.Item(6, 12).Value = 8.5
$ws.Cells.Item(6, 13).Value = 1.07
$ws.Cells.Item(6, 14).Value = 8.5
$ws.Cells.Item(6, 23).Value = 5
$ws.Cells.Item(6, 50).Value = 9

$ws.Cells.Item(7, 7).Value = 1.42
$ws.Cells.Item(7, 8).Value = 4.2
$ws.Cells.Item(7, 9).Value = 8
$ws.Cells.Item(7, 10).Value = 2
$ws.Cells.Item(7, 12).Value = 8
$ws.Cells.Item(7, 13).Value = 1.06
$ws.Cells.Item(7, 14).Value = 10
$ws.Cells.Item(7, 17).Value = 2.07
$ws.Cells.Item(7, 18).Value = 1.69
$ws.Cells.Item(7, 21).Value = 2.38
$ws.Cells.Item(7, 22).Value = 1.53
$ws.Cells.Item(7, 26).Value = 9
$ws.Cells.Item(7, 27).Value = 15
$ws.Cells.Item(7, 31).Value = 23
$ws.Cells.Item(7, 37).Value = 101
$ws.Cells.Item(7, 38).Value = 67
$ws.Cells.Item(7, 40).Value = 3.2
$ws.Cells.Item(7, 41).Value = 7
$ws.Cells.Item(7, 42).Value = 23
$ws.Cells.Item(7, 43).Value = 21
$ws.Cells.Item(7, 50).Value = 8.5
$ws.Cells.Item(7, 52).Value = 51
$ws.Cells.Item(7, 53).Value = 201
$ws.Cells.Item(7, 54).Value = 251

$ws.Cells.Item(8, 13).Value = 1.06
$ws.Cells.Item(8, 14).Value = 10
$ws.Cells.Item(8, 18).Value = 1.63

$ws.Cells.Item(9, 7).Value = 1.7
$ws.Cells.Item(9, 8).Value = 3.5
$ws.Cells.Item(9, 9).Value = 5.5
$ws.Cells.Item(9, 10).Value = 2.38
$ws.Cells.Item(9, 11).Value = 2.05
$ws.Cells.Item(9, 12).Value = 6
$ws.Cells.Item(9, 13).Value = 1.08
$ws.Cells.Item(9, 14).Value = 8
$ws.Cells.Item(9, 17).Value = 2.3
$ws.Cells.Item(9, 18).Value = 1.57
$ws.Cells.Item(9, 21).Value = 2.2
$ws.Cells.Item(9, 22).Value = 1.62
$ws.Cells.Item(9, 24).Value = 7
$ws.Cells.Item(9, 26).Value = 12
$ws.Cells.Item(9, 27).Value = 17
$ws.Cells.Item(9, 30).Value = 7
$ws.Cells.Item(9, 31).Value = 21
$ws.Cells.Item(9, 34).Value = 11
$ws.Cells.Item(9, 35).Value = 26
$ws.Cells.Item(9, 36).Value = 19
$ws.Cells.Item(9, 37).Value = 51
$ws.Cells.Item(9, 38).Value = 41
$ws.Cells.Item(9, 39).Value = 51
$ws.Cells.Item(9, 40).Value = 3.5
$ws.Cells.Item(9, 41).Value = 9
$ws.Cells.Item(9, 43).Value = 34
$ws.Cells.Item(9, 47).Value = 9.5
$ws.Cells.Item(9, 50).Value = 7
$ws.Cells.Item(9, 51).Value = 34
$ws.Cells.Item(9, 53).Value = 126
$ws.Cells.Item(9, 54).Value = 151

$ws.Cells.Item(10, 13).Value = 1.07
$ws.Cells.Item(10, 15).Value = 1.47

$ws.Cells.Item(11, 13).Value = 1.03
$ws.Cells.Item(11, 15).Value = 1.19
$ws.Cells.Item(11, 17).Value = 1.75
$ws.Cells.Item(11, 18).Value = 2.05

$ws.Cells.Item(12, 13).Value = 1.03
$ws.Cells.Item(12, 14).Value = 13
$ws.Cells.Item(12, 15).Value = 1.19
$ws.Cells.Item(12, 17).Value = 1.75
$ws.Cells.Item(12, 18).Value = 2.05

$ws.Cells.Item(13, 21).Value = 1.92
$ws.Cells.Item(13, 22).Value = 1.77

$ws.Cells.Item(14, 11).Value = 2.4
$ws.Cells.Item(14, 19).Value = 1.3
$ws.Cells.Item(14, 20).Value = 3.4
$ws.Cells.Item(14, 21).Value = 1.69
$ws.Cells.Item(14, 29).Value = 15
$ws.Cells.Item(14, 32).Value = 41
$ws.Cells.Item(14, 34).Value = 8.5
$ws.Cells.Item(14, 37).Value = 12
$ws.Cells.Item(14, 39).Value = 21
$ws.Cells.Item(14, 40).Value = 7
$ws.Cells.Item(14, 42).Value = 29
$ws.Cells.Item(14, 46).Value = 3.4
$ws.Cells.Item(14, 55).Value = 101

$ws.Cells.Item(16, 7).Value = 1.75
$ws.Cells.Item(16, 9).Value = 4.5
$ws.Cells.Item(16, 12).Value = 5
$ws.Cells.Item(16, 13).Value = 1.07
$ws.Cells.Item(16, 14).Value = 8.5
$ws.Cells.Item(16, 15).Value = 1.36
$ws.Cells.Item(16, 16).Value = 3
$ws.Cells.Item(16, 17).Value = 2.15
$ws.Cells.Item(16, 18).Value = 1.67
$ws.Cells.Item(16, 21).Value = 2
$ws.Cells.Item(16, 22).Value = 1.73
$ws.Cells.Item(16, 23).Value = 6
$ws.Cells.Item(16, 24).Value = 8
$ws.Cells.Item(16, 25).Value = 9
$ws.Cells.Item(16, 27).Value = 17
$ws.Cells.Item(16, 29).Value = 8.5
$ws.Cells.Item(16, 33).Value = 351
$ws.Cells.Item(16, 37).Value = 51
$ws.Cells.Item(16, 42).Value = 23
$ws.Cells.Item(16, 45).Value = 201
$ws.Cells.Item(16, 47).Value = 9
$ws.Cells.Item(16, 48).Value = 67
$ws.Cells.Item(16, 51).Value = 26
$ws.Cells.Item(16, 54).Value = 126
$ws.Cells.Item(16, 55).Value = 301

$ws.Cells.Item(18, 7).Value = 3.1
$ws.Cells.Item(18, 8).Value = 3
$ws.Cells.Item(18, 9).Value = 2.45
$ws.Cells.Item(18, 10).Value = 3.6
$ws.Cells.Item(18, 11).Value = 2.05
$ws.Cells.Item(18, 12).Value = 3.1
$ws.Cells.Item(18, 17).Value = 2.08
$ws.Cells.Item(18, 18).Value = 1.73
$ws.Cells.Item(18, 21).Value = 1.77
$ws.Cells.Item(18, 22).Value = 1.92
$ws.Cells.Item(18, 27).Value = 26
$ws.Cells.Item(18, 34).Value = 7.5
$ws.Cells.Item(18, 35).Value = 11
$ws.Cells.Item(18, 36).Value = 9.5
$ws.Cells.Item(18, 50).Value = 4.33
$ws.Cells.Item(18, 51).Value = 13
$ws.Cells.Item(18, 53).Value = 41

$ws.Cells.Item(21, 10).Value = 8
$ws.Cells.Item(21, 13).Value = 1.03
$ws.Cells.Item(21, 14).Value = 15
$ws.Cells.Item(21, 22).Value = 1.67
$ws.Cells.Item(21, 25).Value = 23
$ws.Cells.Item(21, 26).Value = 101
$ws.Cells.Item(21, 29).Value = 12
$ws.Cells.Item(21, 33).Value = 451

$ws.Cells.Item(22, 7).Value = 1.85
$ws.Cells.Item(22, 9).Value = 3.6
$ws.Cells.Item(22, 12).Value = 4.33
$ws.Cells.Item(22, 13).Value = 1.02
$ws.Cells.Item(22, 14).Value = 11
$ws.Cells.Item(22, 21).Value = 1.77
$ws.Cells.Item(22, 22).Value = 1.87
$ws.Cells.Item(22, 26).Value = 15
$ws.Cells.Item(22, 33).Value = 600
$ws.Cells.Item(22, 50).Value = 6
$ws.Cells.Item(22, 55).Value = 400

$ws.Cells.Item(23, 7).Value = 2.88
$ws.Cells.Item(23, 12).Value = 2.88
$ws.Cells.Item(23, 21).Value = 1.63

$ws.Cells.Item(24, 7).Value = 1.4
$ws.Cells.Item(24, 10).Value = 1.91
$ws.Cells.Item(24, 11).Value = 2.38
$ws.Cells.Item(24, 13).Value = 1.04
$ws.Cells.Item(24, 14).Value = 9
$ws.Cells.Item(24, 17).Value = 1.7
$ws.Cells.Item(24, 18).Value = 2.1
$ws.Cells.Item(24, 21).Value = 1.87
$ws.Cells.Item(24, 22).Value = 1.77

$ws.Cells.Item(25, 10).Value = 1.65
$ws.Cells.Item(25, 11).Value = 2.6
$ws.Cells.Item(25, 12).Value = 7.6
$ws.Cells.Item(25, 13).Value = 1.01
$ws.Cells.Item(25, 14).Value = 14.6
$ws.Cells.Item(25, 16).Value = 4.9
$ws.Cells.Item(25, 17).Value = 1.52
$ws.Cells.Item(25, 18).Value = 2.22
$ws.Cells.Item(25, 19).Value = 1.27
$ws.Cells.Item(25, 20).Value = 3.52
$ws.Cells.Item(25, 21).Value = 1.93
$ws.Cells.Item(25, 22).Value = 1.7
$ws.Cells.Item(25, 23).Value = 7.8
$ws.Cells.Item(25, 24).Value = 6.5
$ws.Cells.Item(25, 25).Value = 8.75
$ws.Cells.Item(25, 26).Value = 7.7
$ws.Cells.Item(25, 27).Value = 10.5
$ws.Cells.Item(25, 28).Value = 28
$ws.Cells.Item(25, 29).Value = 15
$ws.Cells.Item(25, 31).Value = 23
$ws.Cells.Item(25, 33).Value = 800
$ws.Cells.Item(25, 34).Value = 26
$ws.Cells.Item(25, 38).Value = 120
$ws.Cells.Item(25, 40).Value = 3.15
$ws.Cells.Item(25, 42).Value = 15
$ws.Cells.Item(25, 43).Value = 13
$ws.Cells.Item(25, 44).Value = 37
$ws.Cells.Item(25, 46).Value = 3.35
$ws.Cells.Item(25, 47).Value = 8.75
